# Insert a new weekly price record at row 180 for
# "Terminal La Palmera de La Serena - Papa" (Hortaliza / Fruta y Hortaliza semanal update).
# Inserting a full row shifts every existing row from 180 downward to 181+,
# growing the used range from A1:R276 to A1:R277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(180).Insert()

$newRow = 180

$ws.Cells.Item($newRow,1).Value2  = 8
$ws.Cells.Item($newRow,2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow,3).Value2  = "Coquimbo"
$ws.Cells.Item($newRow,4).Value2  = 44518
$ws.Cells.Item($newRow,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow,5).Value2  = 4
$ws.Cells.Item($newRow,6).Value2  = 100114001
$ws.Cells.Item($newRow,7).Value2  = "Papa"
$ws.Cells.Item($newRow,8).Value2  = "Cardinal"
$ws.Cells.Item($newRow,9).Value2  = "1a nueva(o)"
$ws.Cells.Item($newRow,10).Value2 = 2500
$ws.Cells.Item($newRow,11).Value2 = 11500
$ws.Cells.Item($newRow,12).Value2 = 12000
$ws.Cells.Item($newRow,13).Value2 = 11750
$ws.Cells.Item($newRow,14).Value2 = '$/saco 25 kilos'
$ws.Cells.Item($newRow,15).Value2 = "Provincia del Elquí"
$ws.Cells.Item($newRow,16).Value2 = 470
$ws.Cells.Item($newRow,17).Value2 = 25
$ws.Cells.Item($newRow,18).Value2 = "Hortaliza"
